$d = $word.ActiveDocument

# Set the Normal style (non-heading body text) to Calibri 11pt,
# matching the commit: "Set font to Calibri for non-heading text."
$normal = $d.Styles.Item("Normal")
$normal.Font.Name = "Calibri"
$normal.Font.Size = 11

# Locate "New Computer Language Func" (the first occurrence, in the
# version table's first data row) so we can re-seat the _GoBack
# bookmark there, mirroring where Word leaves it after the edit.
$findRange = $d.Range(0, 0)
$found = $findRange.Find.Execute("New Computer Language Func", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $markRange = $d.Range($findRange.End, $findRange.End)
    $d.Bookmarks.Add("_GoBack", $markRange)
}
